$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Talent_Acquisition (row/col additions + new scenario row)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Talent_Acquisition")

# New data row 7 - it mirrors row 6 (same formats across every column), so
# copy row 6 down (values + styles) and then touch up the handful of cells
# that differ for this scenario.
$ws1.Range("A6:BB6").Copy($ws1.Range("A7:BB7"))
$ws1.Range("A7").Value = "EDIT_PROJECTED_ENDDATE"

# New header cells (BC1:BE1) - give them the same "applied alignment" style
# bucket that the source workbook uses, then fill in the header captions.
$ws1.Range("BC1:BE1").WrapText = $false
$ws1.Range("BC1").Value = "projectedEndDate"
$ws1.Range("BD1").Value = "correctEmploymentAction"
$ws1.Range("BE1").Value = "correctEmploymentActionReason"

$ws1.Range("AR7").Value = "3259988"

# New trailing columns for row 7.
$ws1.Range("AS7").Copy($ws1.Range("BD7"))
$ws1.Range("BD7").Value = "Add Pending Worker"
$ws1.Range("AS7").Copy($ws1.Range("BE7"))
$ws1.Range("BE7").Value = "Additional Hire"
$ws1.Range("BC7").NumberFormat = "d-mmm-yy"
$ws1.Range("BC7").Value = 43687

# Size the new columns to fit their content.
$ws1.Columns.Item(55).ColumnWidth = 9.022135416666666
$ws1.Columns.Item(56).ColumnWidth = 23.736979166666668
$ws1.Columns.Item(57).ColumnWidth = 30.451822916666668

# ---------------------------------------------------------------------------
# Sheet: Line_Manager loses the tab selection / A1:G2 selection it had...
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Line_Manager")
$ws4.Range("D6").Select()

# ...and Talent_Acquisition becomes the active tab / sheet instead, parked
# on the freshly added BE7 cell.
$ws1.Activate()
$ws1.Range("AT1").Select()
$ws1.Range("BE7").Select()
